$d = $word.ActiveDocument

# 1. "B1 – Regression with R" -> "B1 – " stays regular, "Regression with R" becomes bold
$r1 = $d.Content
$r1.Find.Execute("Regression with R", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r1.Font.Bold = 1

# 2. Remove the _GoBack bookmark that currently sits right after "Lunch"
if ($d.Bookmarks.Exists("_GoBack")) {
    $bmOld = $d.Bookmarks.Item("_GoBack")
    $bmOld.Delete()
}

# 3. "C - Supervised Learning: " -> "C -" loses bold, " Supervised Learning: " stays bold
$r3 = $d.Content
$r3.Find.Execute("C - Supervised Learning: Bagging and Boosting, tree-methods", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$cMinus = $d.Range($r3.Start, $r3.Start + 3)
$cMinus.Font.Bold = 0

# 4. Re-insert the _GoBack bookmark between "tre" and "e-methods" inside "tree-methods"
$r4 = $d.Content
$r4.Find.Execute("tree-methods", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPoint = $r4.Start + 3
$insertRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $insertRange)
